$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1443736666666667
$ws.Range("H2").Value = 0.433121
$ws.Range("I2").Value = 0.7378778224885942
$ws.Range("J2").Value = 0.7378778224885942
$ws.Range("M2").Value = 102.3975143333333
$ws.Range("N2").Value = 307.192543
$ws.Range("O2").Value = 0.822753865655704
$ws.Range("P2").Value = 0.8227538656557041
$ws.Range("Q2").Value = 14.78350460185589
$ws.Range("R2").Value = 133.051541416703
$ws.Range("S2").Value = 0.6070918308341042
$ws.Range("T2").Value = 0.6070918308341043

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1443736666666667
$ws.Range("H3").Value = 0.433121
$ws.Range("I3").Value = 0.7378778224885942
$ws.Range("J3").Value = 0.7378778224885942
$ws.Range("O3").Value = 0.1536069072592176
$ws.Range("P3").Value = 0.1536069072592176
$ws.Range("Q3").Value = 2.760058038176112
$ws.Range("R3").Value = 24.840522343585
$ws.Range("S3").Value = 0.1133431302476389
$ws.Range("T3").Value = 0.1133431302476389

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1443736666666667
$ws.Range("H4").Value = 0.433121
$ws.Range("I4").Value = 0.7378778224885942
$ws.Range("J4").Value = 0.7378778224885942
$ws.Range("M4").Value = 2.600356333333333
$ws.Range("N4").Value = 7.801069
$ws.Range("O4").Value = 0.02089360507685526
$ws.Range("P4").Value = 0.02089360507685526
$ws.Range("Q4").Value = 0.3754229784832223
$ws.Range("R4").Value = 3.378806806349
$ws.Range("S4").Value = 0.01541692781804659
$ws.Range("T4").Value = 0.0154169278180466

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1443736666666667
$ws.Range("H5").Value = 0.433121
$ws.Range("I5").Value = 0.7378778224885942
$ws.Range("J5").Value = 0.7378778224885942
$ws.Range("M5").Value = 0.341712
$ws.Range("N5").Value = 1.025136
$ws.Range("O5").Value = 0.002745622008223115
$ws.Range("P5").Value = 0.002745622008223116
$ws.Range("Q5").Value = 0.04933421438400001
$ws.Range("R5").Value = 0.444007929456
$ws.Range("S5").Value = 0.002025933588804434
$ws.Range("T5").Value = 0.002025933588804434

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.051287
$ws.Range("H6").Value = 0.153861
$ws.Range("I6").Value = 0.2621221775114058
$ws.Range("J6").Value = 0.2621221775114058
$ws.Range("M6").Value = 102.3975143333333
$ws.Range("N6").Value = 307.192543
$ws.Range("O6").Value = 0.822753865655704
$ws.Range("P6").Value = 0.8227538656557041
$ws.Range("Q6").Value = 5.251661317613666
$ws.Range("R6").Value = 47.264951858523
$ws.Range("S6").Value = 0.2156620348215997
$ws.Range("T6").Value = 0.2156620348215998

# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.051287
$ws.Range("H7").Value = 0.153861
$ws.Range("I7").Value = 0.2621221775114058
$ws.Range("J7").Value = 0.2621221775114058
$ws.Range("O7").Value = 0.1536069072592176
$ws.Range("P7").Value = 0.1536069072592176
$ws.Range("Q7").Value = 0.9804772564983334
$ws.Range("R7").Value = 8.824295308485
$ws.Range("S7").Value = 0.04026377701157868
$ws.Range("T7").Value = 0.04026377701157869

# Row 8
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.051287
$ws.Range("H8").Value = 0.153861
$ws.Range("I8").Value = 0.2621221775114058
$ws.Range("J8").Value = 0.2621221775114058
$ws.Range("M8").Value = 2.600356333333333
$ws.Range("N8").Value = 7.801069
$ws.Range("O8").Value = 0.02089360507685526
$ws.Range("P8").Value = 0.02089360507685526
$ws.Range("Q8").Value = 0.1333644752676667
$ws.Range("R8").Value = 1.200280277409
$ws.Range("S8").Value = 0.005476677258808662
$ws.Range("T8").Value = 0.005476677258808665

# Row 9
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.051287
$ws.Range("H9").Value = 0.153861
$ws.Range("I9").Value = 0.2621221775114058
$ws.Range("J9").Value = 0.2621221775114058
$ws.Range("M9").Value = 0.341712
$ws.Range("N9").Value = 1.025136
$ws.Range("O9").Value = 0.002745622008223115
$ws.Range("P9").Value = 0.002745622008223116
$ws.Range("Q9").Value = 0.017525383344
$ws.Range("R9").Value = 0.157728450096
$ws.Range("S9").Value = 0.0007196884194186819
$ws.Range("T9").Value = 0.0007196884194186821
